{"js": "const replacements = [\n  [\"335\u00d79=3015\", \"110\u00d77=770\"],\n  [\"122\u00d79=1098\", \"324\u00d72=648\"],\n  [\"339\u00d72=678\", \"187\u00d78=1496\"],\n  [\"953\u00d74=3812\", \"159\u00d77=1113\"],\n  [\"889\u00d76=5334\", \"903\u00d72=1806\"],\n  [\"865\u00d77=6055\", \"578\u00d75=2890\"],\n  [\"470\u00d75=2350\", \"828\u00d76=4968\"],\n  [\"736\u00d75=3680\", \"932\u00d79=8388\"],\n  [\"649\u00d76=3894\", \"502\u00d78=4016\"],\n  [\"128\u00d72=256\", \"231\u00d73=693\"],\n  [\"112\u00d77=784\", \"199\u00d77=1393\"],\n  [\"214\u00d78=1712\", \"852\u00d74=3408\"],\n  [\"239\u00d79=2151\", \"840\u00d76=5040\"],\n  [\"984\u00d78=7872\", \"414\u00d72=828\"],\n  [\"746\u00d79=6714\", \"625\u00d73=1875\"],\n  [\"358\u00d75=1790\", \"596\u00d73=1788\"],\n  [\"821\u00d76=4926\", \"589\u00d79=5301\"],\n  [\"862\u00d74=3448\", \"745\u00d73=2235\"],\n  [\"309\u00d74=1236\", \"715\u00d76=4290\"],\n  [\"490\u00d78=3920\", \"400\u00d72=800\"],\n  [\"390\u00d78=3120\", \"320\u00d72=640\"],\n  [\"140\u00d73=420\", \"827\u00d76=4962\"],\n  [\"539\u00d76=3234\", \"303\u00d77=2121\"],\n  [\"984\u00d74=3936\", \"873\u00d76=5238\"],\n  [\"338\u00d73=1014\", \"269\u00d78=2152\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"335\u00d79=3015\", \"110\u00d77=770\"),\n    @(\"122\u00d79=1098\", \"324\u00d72=648\"),\n    @(\"339\u00d72=678\", \"187\u00d78=1496\"),\n    @(\"953\u00d74=3812\", \"159\u00d77=1113\"),\n    @(\"889\u00d76=5334\", \"903\u00d72=1806\"),\n    @(\"865\u00d77=6055\", \"578\u00d75=2890\"),\n    @(\"470\u00d75=2350\", \"828\u00d76=4968\"),\n    @(\"736\u00d75=3680\", \"932\u00d79=8388\"),\n    @(\"649\u00d76=3894\", \"502\u00d78=4016\"),\n    @(\"128\u00d72=256\", \"231\u00d73=693\"),\n    @(\"112\u00d77=784\", \"199\u00d77=1393\"),\n    @(\"214\u00d78=1712\", \"852\u00d74=3408\"),\n    @(\"239\u00d79=2151\", \"840\u00d76=5040\"),\n    @(\"984\u00d78=7872\", \"414\u00d72=828\"),\n    @(\"746\u00d79=6714\", \"625\u00d73=1875\"),\n    @(\"358\u00d75=1790\", \"596\u00d73=1788\"),\n    @(\"821\u00d76=4926\", \"589\u00d79=5301\"),\n    @(\"862\u00d74=3448\", \"745\u00d73=2235\"),\n    @(\"309\u00d74=1236\", \"715\u00d76=4290\"),\n    @(\"490\u00d78=3920\", \"400\u00d72=800\"),\n    @(\"390\u00d78=3120\", \"320\u00d72=640\"),\n    @(\"140\u00d73=420\", \"827\u00d76=4962\"),\n    @(\"539\u00d76=3234\", \"303\u00d77=2121\"),\n    @(\"984\u00d74=3936\", \"873\u00d76=5238\"),\n    @(\"338\u00d73=1014\", \"269\u00d78=2152\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $find\n    $rng.Find.Replacement.Text = $replace\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n\n$d.Save()\n"}
